$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9456609487533569
$ws.Range("B1").Value = 2.014955759048462
$ws.Range("C1").Value = 7.667108535766602
$ws.Range("D1").Value = 2.675674438476562
$ws.Range("E1").Value = 0.9553760290145874
